$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object 'object[,]' 24,5
$bf[0,0] = 1.02
$bf[0,1] = 1.048995049297712
$bf[0,2] = 1.055427139596091
$bf[0,3] = 1.056177801244151
$bf[0,4] = 1.066690171761096
$bf[1,0] = 1.02
$bf[1,1] = 1.049845166754434
$bf[1,2] = 1.056086856370992
$bf[1,3] = 1.056915886730987
$bf[1,4] = 1.06746585325033
$bf[2,0] = 1.02
$bf[2,1] = 1.050395919622084
$bf[2,2] = 1.056514287087109
$bf[2,3] = 1.057394412396004
$bf[2,4] = 1.067968690458154
$bf[3,0] = 1.02
$bf[3,1] = 1.050627614821045
$bf[3,2] = 1.056694108769313
$bf[3,3] = 1.057595806529154
$bf[3,4] = 1.068180301347833
$bf[4,0] = 1.02
$bf[4,1] = 1.050666526720589
$bf[4,2] = 1.056724309179957
$bf[4,3] = 1.057629634449964
$bf[4,4] = 1.068215844480776
$bf[5,0] = 1.02
$bf[5,1] = 1.050399014923257
$bf[5,2] = 1.056516689365463
$bf[5,3] = 1.057397102564568
$bf[5,4] = 1.06797151715862
$bf[6,0] = 1.02
$bf[6,1] = 1.049282210708367
$bf[6,2] = 1.055649979021987
$bf[6,3] = 1.056427045951694
$bf[6,4] = 1.066952125389382
$bf[7,0] = 1.02
$bf[7,1] = 1.047319466420207
$bf[7,2] = 1.054127016535562
$bf[7,3] = 1.054724927129505
$bf[7,4] = 1.065162956127727
$bf[8,0] = 1.02
$bf[8,1] = 1.046014573491202
$bf[8,2] = 1.053114702773836
$bf[8,3] = 1.05359516556281
$bf[8,4] = 1.063975088296712
$bf[9,0] = 1.02
$bf[9,1] = 1.045450416575334
$bf[9,2] = 1.052677093046873
$bf[9,3] = 1.053107172505134
$bf[9,4] = 1.063461918425263
$bf[10,0] = 1.02
$bf[10,1] = 1.045240996060808
$bf[10,2] = 1.052514656673264
$bf[10,3] = 1.052926092562687
$bf[10,4] = 1.063271484311817
$bf[11,0] = 1.02
$bf[11,1] = 1.045285911449194
$bf[11,2] = 1.052549494757407
$bf[11,3] = 1.052964926539489
$bf[11,4] = 1.063312324896444
$bf[12,0] = 1.02
$bf[12,1] = 1.045433103097669
$bf[12,2] = 1.052663663722983
$bf[12,3] = 1.053092200648303
$bf[12,4] = 1.063446173397824
$bf[13,0] = 1.02
$bf[13,1] = 1.045523810337195
$bf[13,2] = 1.052734021803303
$bf[13,3] = 1.053170642656223
$bf[13,4] = 1.063528665792988
$bf[14,0] = 1.02
$bf[14,1] = 1.046052033211829
$bf[14,2] = 1.053143761004902
$bf[14,3] = 1.053627577528097
$bf[14,4] = 1.064009170825574
$bf[15,0] = 1.02
$bf[15,1] = 1.046383607695574
$bf[15,2] = 1.053400976023854
$bf[15,3] = 1.053914523659798
$bf[15,4] = 1.064310897560919
$bf[16,0] = 1.02
$bf[16,1] = 1.04657709334082
$bf[16,2] = 1.053551075387346
$bf[16,3] = 1.054082010247355
$bf[16,4] = 1.064487003834696
$bf[17,0] = 1.02
$bf[17,1] = 1.046643081117174
$bf[17,2] = 1.053602267240758
$bf[17,3] = 1.054139138421737
$bf[17,4] = 1.064547070827076
$bf[18,0] = 1.02
$bf[18,1] = 1.046348024205113
$bf[18,2] = 1.05337337201254
$bf[18,3] = 1.053883725067472
$bf[18,4] = 1.064278513313532
$bf[19,0] = 1.02
$bf[19,1] = 1.045389755154795
$bf[19,2] = 1.052630040719041
$bf[19,3] = 1.053054716544529
$bf[19,4] = 1.063406753364353
$bf[20,0] = 1.02
$bf[20,1] = 1.044788020732125
$bf[20,2] = 1.052163324070121
$bf[20,3] = 1.05253454229614
$bf[20,4] = 1.062859685862047
$bf[21,0] = 1.02
$bf[21,1] = 1.045106938220668
$bf[21,2] = 1.052410677612115
$bf[21,3] = 1.052810195731375
$bf[21,4] = 1.063149597193443
$bf[22,0] = 1.02
$bf[22,1] = 1.046364102593879
$bf[22,2] = 1.053385844860406
$bf[22,3] = 1.053897641268198
$bf[22,4] = 1.06429314600845
$bf[23,0] = 1.02
$bf[23,1] = 1.047826255251105
$bf[23,2] = 1.054520219000314
$bf[23,3] = 1.055164095411742
$bf[23,4] = 1.065624642210949
$ws.Range("B2:F25").Value = $bf

$im = New-Object 'object[,]' 24,5
$im[0,0] = 1.049910084595803
$im[0,1] = 1.054035845434295
$im[0,2] = 1.058167627690304
$im[0,3] = 1.058916226787733
$im[0,4] = 1.069400043080338
$im[1,0] = 1.050145554399193
$im[1,1] = 1.054535551589985
$im[1,2] = 1.05864117586472
$im[1,3] = 1.05946809219971
$im[1,4] = 1.069991465215183
$im[2,0] = 1.050296967654735
$im[2,1] = 1.05485888137687
$im[2,2] = 1.058947456738365
$im[2,3] = 1.059825450192251
$im[2,4] = 1.07037441003271
$im[3,0] = 1.050360393150236
$im[3,1] = 1.054994804642434
$im[3,2] = 1.0590761831881
$im[3,3] = 1.059975745253874
$im[3,4] = 1.070535459514143
$im[4,0] = 1.050371029137728
$im[4,1] = 1.055017626449398
$im[4,2] = 1.059097794890363
$im[4,3] = 1.060000984034708
$im[4,4] = 1.070562503872052
$im[5,0] = 1.050297816048768
$im[5,1] = 1.054860697609865
$im[5,2] = 1.058949176922906
$im[5,3] = 1.0598274582014
$im[5,4] = 1.070376561752061
$im[6,0] = 1.04998985920692
$im[6,1] = 1.054204725459441
$im[6,2] = 1.058327692950028
$im[6,3] = 1.059102676699378
$im[6,4] = 1.069599863165911
$im[7,0] = 1.049439950881012
$im[7,1] = 1.053048774617575
$im[7,2] = 1.057231573827144
$im[7,3] = 1.057827607789805
$im[7,4] = 1.068233246038716
$im[8,0] = 1.049068522690467
$im[8,1] = 1.052278191139925
$im[8,2] = 1.05650024668469
$im[8,3] = 1.056979054540224
$im[8,4] = 1.067323625121014
$im[9,0] = 1.048906557100461
$im[9,1] = 1.051944548836872
$im[9,2] = 1.05618345327096
$im[9,3] = 1.056611994326801
$im[9,4] = 1.066930115050192
$im[10,0] = 1.048846226241561
$im[10,1] = 1.051820624408914
$im[10,2] = 1.056065764850079
$im[10,3] = 1.05647570875754
$im[10,4] = 1.066784003875914
$im[11,0] = 1.048859175077863
$im[11,1] = 1.051847206363216
$im[11,2] = 1.056091010173965
$im[11,3] = 1.056504939874042
$im[11,4] = 1.06681534265905
$im[12,0] = 1.048901573588246
$im[12,1] = 1.051934305091495
$im[12,2] = 1.056173725453136
$im[12,3] = 1.056600727748909
$im[12,4] = 1.066918036303391
$im[13,0] = 1.048927674258124
$im[13,1] = 1.051987970219492
$im[13,2] = 1.056224686823092
$im[13,3] = 1.056659753409854
$im[13,4] = 1.066981316714628
$im[14,0] = 1.049079247978556
$im[14,1] = 1.052300334532454
$im[14,2] = 1.056521268726691
$im[14,3] = 1.057003423019638
$im[14,4] = 1.067349748850254
$im[15,0] = 1.049174022970679
$im[15,1] = 1.052496280180198
$im[15,2] = 1.056707274511175
$im[15,3] = 1.057219097727448
$im[15,4] = 1.067580954608863
$im[16,0] = 1.049229194077557
$im[16,1] = 1.052610574323691
$im[16,2] = 1.056815756487959
$im[16,3] = 1.057344932632181
$im[16,4] = 1.067715847814145
$im[17,0] = 1.049247987392996
$im[17,1] = 1.052649546032354
$im[17,2] = 1.056852743952744
$im[17,3] = 1.057387845057786
$im[17,4] = 1.067761848763863
$im[18,0] = 1.049163865835312
$im[18,1] = 1.05247525681246
$im[18,2] = 1.056687319087322
$im[18,3] = 1.057195954189968
$im[18,4] = 1.067556144807979
$im[19,0] = 1.048889092959373
$im[19,1] = 1.051908656533025
$im[19,2] = 1.056149368329938
$im[19,3] = 1.056572519023357
$im[19,4] = 1.066887794026083
$im[20,0] = 1.048715351359387
$im[20,1] = 1.051552443084687
$im[20,2] = 1.055811038263826
$im[20,3] = 1.056180870757872
$im[20,4] = 1.066467900094047
$im[21,0] = 1.048807547769623
$im[21,1] = 1.051741275215845
$im[21,2] = 1.055990402364243
$im[21,3] = 1.056388459114254
$im[21,4] = 1.066690462486442
$im[22,0] = 1.049168455746148
$im[22,1] = 1.052484756360272
$im[22,2] = 1.056696336121797
$im[22,3] = 1.057206411648868
$im[22,4] = 1.067567355181762
$im[23,0] = 1.049582969142592
$im[23,1] = 1.053347612594776
$im[23,2] = 1.057515055252496
$im[23,3] = 1.058156987036892
$im[23,4] = 1.068586299489538
$ws.Range("I2:M25").Value = $im

Write-Output "vm_pu values updated"